# Fruta / hortaliza, semanal
# Update the "Fecha" (D), "Volumen" (M), "Precio mínimo" (N), "Precio máximo" (O),
# "Precio promedio ponderado" (P) and "Precio $/Kg" (S) columns for rows 2-13.
# This reflects a re-sort / re-shuffle of the weekly price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio, Precio $/Kg)
$rows = @{
    2  = @(44186, 40, 15000, 15000, 15000, 3000)
    3  = @(44902, 35, 12000, 12000, 12000, 2400)
    4  = @(44193, 40, 15000, 15000, 15000, 3000)
    5  = @(44907, 45, 25000, 25000, 25000, 5000)
    6  = @(44196, 56, 15000, 15000, 15000, 3000)
    7  = @(44189, 40, 15000, 15000, 15000, 3000)
    8  = @(44931, 50, 18000, 18000, 18000, 3600)
    9  = @(44188, 30, 15000, 15000, 15000, 3000)
    10 = @(44914, 56, 23000, 23000, 23000, 4600)
    11 = @(44179, 45, 20000, 20000, 20000, 4000)
    12 = @(44175, 25, 20000, 20000, 20000, 4000)
    13 = @(44181, 30, 20000, 20000, 20000, 4000)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 13).Value = $vals[1]   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals[2]   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals[3]   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals[4]   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $vals[5]   # S - Precio $/Kg
}
